$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @{
    2 = @(922.356803894043, 29.77652573723383, 29.52512967424755)
    3 = @(1800.029103040695, 172.0132764561589, 61.94300509923724)
    4 = @(1800.024947881699, 53.61333600293771, 37.9003799541165)
    5 = @(1001.752865076065, 97.66105119990154, 97.38106729892448)
    6 = @(1800.015377044678, 78.94737691767912, 50.74460587833157)
    7 = @(1800.007972002029, 63.27065289641369, 40.83998661785475)
    8 = @(1056.470659971237, 75.16669942679187, 74.47514465578878)
    9 = @(1800.006425142288, 93.47627693142671, 48.04548705034094)
    10 = @(1657.981845140457, 62.13426192009212, 61.66116771814921)
    11 = @(1800.014371871948, 16.91760078202191, 1.876862522905473)
    12 = @(1800.008924007416, 54.51425737646211, 40.11320110412142)
    13 = @(1800.025357961655, 109.8504091104679, 58.31575468341022)
    14 = @(1800.035519123077, 58.47945732726614, 19.39498828940649)
    15 = @(18.12905406951904, 18.11557489757924, 18.11557489757924)
    16 = @(1800.016659021378, 71.75748607984156, 28.78186443296211)
    17 = @(1800.007781982422, 31.06558919058114, 15.35776508589702)
    18 = @(1800.012269020081, 106.5975133184138, 74.6320076280808)
    19 = @(1800.019261837006, 25.630824372315, 0)
    20 = @(1800.02040719986, 84.21905247243454, 62.99160145332882)
    21 = @(1800.013839960098, 31.38261863492477, 0.08614865408327427)
    22 = @(1800.021936893463, 35.46502725731479, 0)
    23 = @(1800.019460916519, 81.07479082057147, 43.03001891869729)
    24 = @(1800.025895833969, 94.85758957152964, 49.12991088775048)
    25 = @(1800.017073869705, 93.48909923890439, 41.18696053878347)
    26 = @(1800.017382144928, 163.9716926168116, 73.65684718622603)
    27 = @(1800.011164903641, 74.42324030702281, 28.36626460955025)
    28 = @(1800.013070106506, 54.46917897083244, 14.04922060158905)
    29 = @(1800.044002056122, 51.63601574645354, 14.46231382738935)
    30 = @(1800.037311792374, 95.9879650954289, 27.69666180211224)
    31 = @(931.1446940898895, 35.00571185716035, 34.85394514358636)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
}
